$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stale "Source" date-stamp row and the extra blank row ---
# Old layout: row1=Source/path, row2=date, row3=blank, row4=blank, row5=Component Count, row6=headers, row7..=data
# New layout: row1=Source/name, row2=blank, row3=Component Count, row4=headers, row5..=data
$ws.Rows("2:3").Delete()

# --- Replace the old source-file date/path value with the schematic name ---
$ws.Range("B1").Value = "AutoDischarger"
$ws.Range("B1").HorizontalAlignment = 1

# --- Add a pullup resistor (R49, 10K) to the BOM, right before the switch (S1) row ---
$ws.Rows("20:20").Insert()
$ws.Range("A20").Value = "R49"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "10K"
$ws.Range("D20").Value = "Resistor_SMD:R_0603_1608Metric"

# --- Minor column B width tweak ---
$ws.Columns("B:B").ColumnWidth = 8.5

# --- Update the active selection to reflect where the edit was made ---
$ws.Range("C3").Select()
